# Natmi following Dr Hou advice:
# Refresh the Lgi3-Adam22 sending/target cluster table with updated
# ligand/receptor statistics, expanding the original 4 data rows (FAPs -> *)
# into 8 rows that also cover the new "sCs" sending cluster (FAPs -> *, sCs -> *).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$arr = New-Object 'object[,]' 8,20
$arr[0,0] = "FAPs"
$arr[0,1] = "Lgi3"
$arr[0,2] = "Adam22"
$arr[0,3] = "ECs"
$arr[0,4] = 3
$arr[0,5] = 1
$arr[0,6] = 0.6167776666666667
$arr[0,7] = 1.850333
$arr[0,8] = 0.7836323164322263
$arr[0,9] = 0.7836323164322262
$arr[0,10] = 3
$arr[0,11] = 1
$arr[0,12] = 5.293357333333334
$arr[0,13] = 15.880072
$arr[0,14] = 0.2811762939216672
$arr[0,15] = 0.2811762939216672
$arr[0,16] = 3.264824584886223
$arr[0,17] = 29.383421263976
$arr[0,18] = 0.2203388305316645
$arr[0,19] = 0.2203388305316645
$arr[1,0] = "FAPs"
$arr[1,1] = "Lgi3"
$arr[1,2] = "Adam22"
$arr[1,3] = "FAPs"
$arr[1,4] = 3
$arr[1,5] = 1
$arr[1,6] = 0.6167776666666667
$arr[1,7] = 1.850333
$arr[1,8] = 0.7836323164322263
$arr[1,9] = 0.7836323164322262
$arr[1,10] = 3
$arr[1,11] = 1
$arr[1,12] = 3.580253333333333
$arr[1,13] = 10.74076
$arr[1,14] = 0.1901784255576477
$arr[1,15] = 0.1901784255576477
$arr[1,16] = 2.208220297008888
$arr[1,17] = 19.87398267308
$arr[1,18] = 0.1490299601551732
$arr[1,19] = 0.1490299601551731
$arr[2,0] = "FAPs"
$arr[2,1] = "Lgi3"
$arr[2,2] = "Adam22"
$arr[2,3] = "M2"
$arr[2,4] = 3
$arr[2,5] = 1
$arr[2,6] = 0.6167776666666667
$arr[2,7] = 1.850333
$arr[2,8] = 0.7836323164322263
$arr[2,9] = 0.7836323164322262
$arr[2,10] = 3
$arr[2,11] = 1
$arr[2,12] = 4.115194333333333
$arr[2,13] = 12.345583
$arr[2,14] = 0.2185937994640287
$arr[2,15] = 0.2185937994640287
$arr[2,16] = 2.538159958793222
$arr[2,17] = 22.843439629139
$arr[2,18] = 0.1712971654317184
$arr[2,19] = 0.1712971654317183
$arr[3,0] = "FAPs"
$arr[3,1] = "Lgi3"
$arr[3,2] = "Adam22"
$arr[3,3] = "sCs"
$arr[3,4] = 3
$arr[3,5] = 1
$arr[3,6] = 0.6167776666666667
$arr[3,7] = 1.850333
$arr[3,8] = 0.7836323164322263
$arr[3,9] = 0.7836323164322262
$arr[3,10] = 3
$arr[3,11] = 1
$arr[3,12] = 5.836954666666667
$arr[3,13] = 17.510864
$arr[3,14] = 0.3100514810566565
$arr[3,15] = 0.3100514810566565
$arr[3,16] = 3.600103279745778
$arr[3,17] = 32.400929517712
$arr[3,18] = 0.2429663603136702
$arr[3,19] = 0.2429663603136702
$arr[4,0] = "sCs"
$arr[4,1] = "Lgi3"
$arr[4,2] = "Adam22"
$arr[4,3] = "ECs"
$arr[4,4] = 2
$arr[4,5] = 0.6666666666666666
$arr[4,6] = 0.1702976666666667
$arr[4,7] = 0.510893
$arr[4,8] = 0.2163676835677737
$arr[4,9] = 0.2163676835677737
$arr[4,10] = 3
$arr[4,11] = 1
$arr[4,12] = 5.293357333333334
$arr[4,13] = 15.880072
$arr[4,14] = 0.2811762939216672
$arr[4,15] = 0.2811762939216672
$arr[4,16] = 0.9014464026995557
$arr[4,17] = 8.113017624296
$arr[4,18] = 0.06083746339000261
$arr[4,19] = 0.06083746339000261
$arr[5,0] = "sCs"
$arr[5,1] = "Lgi3"
$arr[5,2] = "Adam22"
$arr[5,3] = "FAPs"
$arr[5,4] = 2
$arr[5,5] = 0.6666666666666666
$arr[5,6] = 0.1702976666666667
$arr[5,7] = 0.510893
$arr[5,8] = 0.2163676835677737
$arr[5,9] = 0.2163676835677737
$arr[5,10] = 3
$arr[5,11] = 1
$arr[5,12] = 3.580253333333333
$arr[5,13] = 10.74076
$arr[5,14] = 0.1901784255576477
$arr[5,15] = 0.1901784255576477
$arr[5,16] = 0.6097087887422222
$arr[5,17] = 5.487379098679999
$arr[5,18] = 0.04114846540247452
$arr[5,19] = 0.04114846540247452
$arr[6,0] = "sCs"
$arr[6,1] = "Lgi3"
$arr[6,2] = "Adam22"
$arr[6,3] = "M2"
$arr[6,4] = 2
$arr[6,5] = 0.6666666666666666
$arr[6,6] = 0.1702976666666667
$arr[6,7] = 0.510893
$arr[6,8] = 0.2163676835677737
$arr[6,9] = 0.2163676835677737
$arr[6,10] = 3
$arr[6,11] = 1
$arr[6,12] = 4.115194333333333
$arr[6,13] = 12.345583
$arr[6,14] = 0.2185937994640287
$arr[6,15] = 0.2185937994640287
$arr[6,16] = 0.7008079928465556
$arr[6,17] = 6.307271935619
$arr[6,18] = 0.04729663403231035
$arr[6,19] = 0.04729663403231035
$arr[7,0] = "sCs"
$arr[7,1] = "Lgi3"
$arr[7,2] = "Adam22"
$arr[7,3] = "sCs"
$arr[7,4] = 2
$arr[7,5] = 0.6666666666666666
$arr[7,6] = 0.1702976666666667
$arr[7,7] = 0.510893
$arr[7,8] = 0.2163676835677737
$arr[7,9] = 0.2163676835677737
$arr[7,10] = 3
$arr[7,11] = 1
$arr[7,12] = 5.836954666666667
$arr[7,13] = 17.510864
$arr[7,14] = 0.3100514810566565
$arr[7,15] = 0.3100514810566565
$arr[7,16] = 0.9940197601724446
$arr[7,17] = 8.946177841552002
$arr[7,18] = 0.06708512074298623
$arr[7,19] = 0.06708512074298623

$ws.Range("A2:T9").Value = $arr